$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.750601
$ws.Range("H2").Value = 65.251803
$ws.Range("I2").Value = 0.3612485837751334
$ws.Range("J2").Value = 0.3773020033645041
$ws.Range("M2").Value = 0.1847186666666667
$ws.Range("N2").Value = 0.554156
$ws.Range("O2").Value = 0.003664790727240103
$ws.Range("P2").Value = 0.00367452998950282
$ws.Range("Q2").Value = 4.017742015918667
$ws.Range("R2").Value = 36.159678143268
$ws.Range("S2").Value = 0.001323900460047728
$ws.Range("T2").Value = 0.001386407526462364
$ws.Range("G3").Value = 21.750601
$ws.Range("H3").Value = 65.251803
$ws.Range("I3").Value = 0.3612485837751334
$ws.Range("J3").Value = 0.3773020033645041
$ws.Range("O3").Value = 0.008410297789740796
$ws.Range("P3").Value = 0.008432648341785422
$ws.Range("Q3").Value = 9.220282769509332
$ws.Range("R3").Value = 82.98254492558399
$ws.Range("S3").Value = 0.003038208165670997
$ws.Range("T3").Value = 0.003181655113024003
$ws.Range("G4").Value = 21.750601
$ws.Range("H4").Value = 65.251803
$ws.Range("I4").Value = 0.3612485837751334
$ws.Range("J4").Value = 0.3773020033645041
$ws.Range("M4").Value = 23.33139033333333
$ws.Range("N4").Value = 69.99417099999999
$ws.Range("O4").Value = 0.4628912956670289
$ws.Range("P4").Value = 0.4641214395041984
$ws.Range("Q4").Value = 507.4717619155903
$ws.Range("R4").Value = 4567.245857240312
$ws.Range("S4").Value = 0.1672188250015507
$ws.Range("T4").Value = 0.1751139489293515
$ws.Range("G5").Value = 21.750601
$ws.Range("H5").Value = 65.251803
$ws.Range("I5").Value = 0.3612485837751334
$ws.Range("J5").Value = 0.3773020033645041
$ws.Range("M5").Value = 0.400781
$ws.Range("N5").Value = 0.801562
$ws.Range("O5").Value = 0.007951435114592365
$ws.Range("P5").Value = 0.005315044152631857
$ws.Range("Q5").Value = 8.717227619380999
$ws.Range("R5").Value = 52.303365716286
$ws.Range("S5").Value = 0.002872444674126358
$ws.Range("T5").Value = 0.002005376806758793
$ws.Range("G6").Value = 21.750601
$ws.Range("H6").Value = 65.251803
$ws.Range("I6").Value = 0.3612485837751334
$ws.Range("J6").Value = 0.3773020033645041
$ws.Range("M6").Value = 26.06280633333334
$ws.Range("N6").Value = 78.18841900000001
$ws.Range("O6").Value = 0.517082180701398
$ws.Range("P6").Value = 0.5184563380118814
$ws.Range("Q6").Value = 566.8817014966064
$ws.Range("R6").Value = 5101.935313469457
$ws.Range("S6").Value = 0.1867952054737377
$ws.Range("T6").Value = 0.1956146149889073
$ws.Range("I7").Value = 0.2797740820980411
$ws.Range("J7").Value = 0.2922068802649305
$ws.Range("M7").Value = 0.1847186666666667
$ws.Range("N7").Value = 0.554156
$ws.Range("O7").Value = 0.003664790727240103
$ws.Range("P7").Value = 0.00367452998950282
$ws.Range("Q7").Value = 3.111597207838667
$ws.Range("R7").Value = 28.004374870548
$ws.Range("S7").Value = 0.001025313461795012
$ws.Range("T7").Value = 0.001073722944672547
$ws.Range("I8").Value = 0.2797740820980411
$ws.Range("J8").Value = 0.2922068802649305
$ws.Range("O8").Value = 0.008410297789740796
$ws.Range("P8").Value = 0.008432648341785422
$ws.Range("S8").Value = 0.002352983344295915
$ws.Range("T8").Value = 0.002464077864324357
$ws.Range("I9").Value = 0.2797740820980411
$ws.Range("J9").Value = 0.2922068802649305
$ws.Range("M9").Value = 23.33139033333333
$ws.Range("N9").Value = 69.99417099999999
$ws.Range("O9").Value = 0.4628912956670289
$ws.Range("P9").Value = 0.4641214395041984
$ws.Range("Q9").Value = 393.0186933798103
$ws.Range("R9").Value = 3537.168240418293
$ws.Range("S9").Value = 0.1295049873564159
$ws.Range("T9").Value = 0.1356194779015905
$ws.Range("I10").Value = 0.2797740820980411
$ws.Range("J10").Value = 0.2922068802649305
$ws.Range("M10").Value = 0.400781
$ws.Range("N10").Value = 0.801562
$ws.Range("O10").Value = 0.007951435114592365
$ws.Range("P10").Value = 0.005315044152631857
$ws.Range("Q10").Value = 6.751180392641
$ws.Range("R10").Value = 40.507082355846
$ws.Range("S10").Value = 0.002224605460547211
$ws.Range("T10").Value = 0.001553092470310916
$ws.Range("I11").Value = 0.2797740820980411
$ws.Range("J11").Value = 0.2922068802649305
$ws.Range("M11").Value = 26.06280633333334
$ws.Range("N11").Value = 78.18841900000001
$ws.Range("O11").Value = 0.517082180701398
$ws.Range("P11").Value = 0.5184563380118814
$ws.Range("Q11").Value = 439.0295625161864
$ws.Range("R11").Value = 3951.266062645678
$ws.Range("S11").Value = 0.144666192474987
$ws.Range("T11").Value = 0.1514965090840321
$ws.Range("G12").Value = 8.938416999999999
$ws.Range("H12").Value = 26.815251
$ws.Range("I12").Value = 0.1484552303838214
$ws.Range("J12").Value = 0.1550523887136425
$ws.Range("M12").Value = 0.1847186666666667
$ws.Range("N12").Value = 0.554156
$ws.Range("O12").Value = 0.003664790727240103
$ws.Range("P12").Value = 0.00367452998950282
$ws.Range("Q12").Value = 1.651092470350667
$ws.Range("R12").Value = 14.859832233156
$ws.Range("S12").Value = 0.0005440573517209218
$ws.Range("T12").Value = 0.0005697446522723278
$ws.Range("G13").Value = 8.938416999999999
$ws.Range("H13").Value = 26.815251
$ws.Range("I13").Value = 0.1484552303838214
$ws.Range("J13").Value = 0.1550523887136425
$ws.Range("O13").Value = 0.008410297789740796
$ws.Range("P13").Value = 0.008432648341785422
$ws.Range("Q13").Value = 3.789078391525333
$ws.Range("R13").Value = 34.101705523728
$ws.Range("S13").Value = 0.001248552695972514
$ws.Range("T13").Value = 0.001307502268575966
$ws.Range("G14").Value = 8.938416999999999
$ws.Range("H14").Value = 26.815251
$ws.Range("I14").Value = 0.1484552303838214
$ws.Range("J14").Value = 0.1550523887136425
$ws.Range("M14").Value = 23.33139033333333
$ws.Range("N14").Value = 69.99417099999999
$ws.Range("O14").Value = 0.4628912956670289
$ws.Range("P14").Value = 0.4641214395041984
$ws.Range("Q14").Value = 208.5456959891023
$ws.Range("R14").Value = 1876.911263901921
$ws.Range("S14").Value = 0.06871863394091438
$ws.Range("T14").Value = 0.07196313784834026
$ws.Range("G15").Value = 8.938416999999999
$ws.Range("H15").Value = 26.815251
$ws.Range("I15").Value = 0.1484552303838214
$ws.Range("J15").Value = 0.1550523887136425
$ws.Range("M15").Value = 0.400781
$ws.Range("N15").Value = 0.801562
$ws.Range("O15").Value = 0.007951435114592365
$ws.Range("P15").Value = 0.005315044152631857
$ws.Range("Q15").Value = 3.582347703677
$ws.Range("R15").Value = 21.494086222062
$ws.Range("S15").Value = 0.001180432131818817
$ws.Range("T15").Value = 0.0008241102919840472
$ws.Range("G16").Value = 8.938416999999999
$ws.Range("H16").Value = 26.815251
$ws.Range("I16").Value = 0.1484552303838214
$ws.Range("J16").Value = 0.1550523887136425
$ws.Range("M16").Value = 26.06280633333334
$ws.Range("N16").Value = 78.18841900000001
$ws.Range("O16").Value = 0.517082180701398
$ws.Range("P16").Value = 0.5184563380118814
$ws.Range("Q16").Value = 232.9602311975744
$ws.Range("R16").Value = 2096.642080778169
$ws.Range("S16").Value = 0.07676355426339482
$ws.Range("T16").Value = 0.08038789365246984
$ws.Range("G17").Value = 7.6853705
$ws.Range("H17").Value = 15.370741
$ws.Range("I17").Value = 0.1276437928732263
$ws.Range("J17").Value = 0.08887741190073968
$ws.Range("M17").Value = 0.1847186666666667
$ws.Range("N17").Value = 0.554156
$ws.Range("O17").Value = 0.003664790727240103
$ws.Range("P17").Value = 0.00367452998950282
$ws.Range("Q17").Value = 1.419631391599333
$ws.Range("R17").Value = 8.517788349596
$ws.Range("S17").Value = 0.0004677877885115561
$ws.Range("T17").Value = 0.0003265827154186627
$ws.Range("G18").Value = 7.6853705
$ws.Range("H18").Value = 15.370741
$ws.Range("I18").Value = 0.1276437928732263
$ws.Range("J18").Value = 0.08887741190073968
$ws.Range("O18").Value = 0.008410297789740796
$ws.Range("P18").Value = 0.008432648341785422
$ws.Range("Q18").Value = 3.257900285074667
$ws.Range("R18").Value = 19.547401710448
$ws.Range("S18").Value = 0.001073522309075827
$ws.Range("T18").Value = 0.0007494719600869524
$ws.Range("G19").Value = 7.6853705
$ws.Range("H19").Value = 15.370741
$ws.Range("I19").Value = 0.1276437928732263
$ws.Range("J19").Value = 0.08887741190073968
$ws.Range("M19").Value = 23.33139033333333
$ws.Range("N19").Value = 69.99417099999999
$ws.Range("O19").Value = 0.4628912956670289
$ws.Range("P19").Value = 0.4641214395041984
$ws.Range("Q19").Value = 179.3103789917852
$ws.Range("R19").Value = 1075.862273950711
$ws.Range("S19").Value = 0.05908520066694161
$ws.Range("T19").Value = 0.04124991235077887
$ws.Range("G20").Value = 7.6853705
$ws.Range("H20").Value = 15.370741
$ws.Range("I20").Value = 0.1276437928732263
$ws.Range("J20").Value = 0.08887741190073968
$ws.Range("M20").Value = 0.400781
$ws.Range("N20").Value = 0.801562
$ws.Range("O20").Value = 0.007951435114592365
$ws.Range("P20").Value = 0.005315044152631857
$ws.Range("Q20").Value = 3.0801504743605
$ws.Range("R20").Value = 12.320601897442
$ws.Range("S20").Value = 0.001014951336811926
$ws.Range("T20").Value = 0.0004723873684240795
$ws.Range("G21").Value = 7.6853705
$ws.Range("H21").Value = 15.370741
$ws.Range("I21").Value = 0.1276437928732263
$ws.Range("J21").Value = 0.08887741190073968
$ws.Range("M21").Value = 26.06280633333334
$ws.Range("N21").Value = 78.18841900000001
$ws.Range("O21").Value = 0.517082180701398
$ws.Range("P21").Value = 0.5184563380118814
$ws.Range("Q21").Value = 200.3023229414132
$ws.Range("R21").Value = 1201.813937648479
$ws.Range("S21").Value = 0.06600233077188543
$ws.Range("T21").Value = 0.04607905750603111
$ws.Range("G22").Value = 4.990062666666667
$ws.Range("H22").Value = 14.970188
$ws.Range("I22").Value = 0.08287831086977776
$ws.Range("J22").Value = 0.08656131575618316
$ws.Range("M22").Value = 0.1847186666666667
$ws.Range("N22").Value = 0.554156
$ws.Range("O22").Value = 0.003664790727240103
$ws.Range("P22").Value = 0.00367452998950282
$ws.Range("Q22").Value = 0.9217577223697778
$ws.Range("R22").Value = 8.295819501327999
$ws.Range("S22").Value = 0.0003037316651648841
$ws.Range("T22").Value = 0.000318072150676918
$ws.Range("G23").Value = 4.990062666666667
$ws.Range("H23").Value = 14.970188
$ws.Range("I23").Value = 0.08287831086977776
$ws.Range("J23").Value = 0.08656131575618316
$ws.Range("O23").Value = 0.008410297789740796
$ws.Range("P23").Value = 0.008432648341785422
$ws.Range("Q23").Value = 2.115334138318222
$ws.Range("R23").Value = 19.038007244864
$ws.Range("S23").Value = 0.0006970312747255425
$ws.Range("T23").Value = 0.0007299411357741423
$ws.Range("G24").Value = 4.990062666666667
$ws.Range("H24").Value = 14.970188
$ws.Range("I24").Value = 0.08287831086977776
$ws.Range("J24").Value = 0.08656131575618316
$ws.Range("M24").Value = 23.33139033333333
$ws.Range("N24").Value = 69.99417099999999
$ws.Range("O24").Value = 0.4628912956670289
$ws.Range("P24").Value = 0.4641214395041984
$ws.Range("Q24").Value = 116.4250998637942
$ws.Range("R24").Value = 1047.825898774148
$ws.Range("S24").Value = 0.03836364870120623
$ws.Range("T24").Value = 0.04017496247413718
$ws.Range("G25").Value = 4.990062666666667
$ws.Range("H25").Value = 14.970188
$ws.Range("I25").Value = 0.08287831086977776
$ws.Range("J25").Value = 0.08656131575618316
$ws.Range("M25").Value = 0.400781
$ws.Range("N25").Value = 0.801562
$ws.Range("O25").Value = 0.007951435114592365
$ws.Range("P25").Value = 0.005315044152631857
$ws.Range("Q25").Value = 1.999922305609333
$ws.Range("R25").Value = 11.999533833656
$ws.Range("S25").Value = 0.000659001511288053
$ws.Range("T25").Value = 0.0004600772151540212
$ws.Range("G26").Value = 4.990062666666667
$ws.Range("H26").Value = 14.970188
$ws.Range("I26").Value = 0.08287831086977776
$ws.Range("J26").Value = 0.08656131575618316
$ws.Range("M26").Value = 26.06280633333334
$ws.Range("N26").Value = 78.18841900000001
$ws.Range("O26").Value = 0.517082180701398
$ws.Range("P26").Value = 0.5184563380118814
$ws.Range("Q26").Value = 130.0550368725303
$ws.Range("R26").Value = 1170.495331852772
$ws.Range("S26").Value = 0.04285489771739306
$ws.Range("T26").Value = 0.04487826278044089
